$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the blank row 2 so the data rows shift up (rows 3-5 -> 2-4)
$ws.Rows.Item(2).Delete()
